# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# For cells whose new text looks like a plain number (e.g. "1.00", "2.70"),
# force the cell to Text format before/after the write so the literal
# string (incl. trailing zeros) round-trips instead of being coerced into
# a numeric cell, then restore the "Normal" style so no stray formatting
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.575.55"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.307.31"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.305.68"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "3.852.29"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "3.308.12"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").Value = "63.659.12"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.735"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  +3.69%  "
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "3.130.21"
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "428.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("E42").Value = "  +8.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.12%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.52%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("E51").Value = "  -0.17%  "
